$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target data for rows 2-32 (Location, Site, Exposure period, Notes)
$data = @(
    @("30/12/20 10:45am- 12:15pm", "", "", ""),
    @("Bentleigh", "Il Centro Deli  5/284/292 Centre Rd, Bentleigh VIC 3204", "22/12/20 12:00pm-12:30pm", "Case shopped in store"),
    @("Black Rock", "Woolworths Metro  40 Bluff Road, Black Rock VIC 3193", "30/12/20 5:30pm-5:55pm", "Case shopped"),
    @("Box Hill South", "Bunnings  259 Middleborough Road, Box Hill South VIC 3128", "30/12/20 12:00pm-12:40pm", "Case shopped"),
    @("Brighton", "Brighton Beach  Brighton, VIC 3186", "29/12/20 12:00pm-3:00pm", "Case attended beach"),
    @("Camberwell", "Fu Lin Asian Grocery Supermarket  1397 Toorak Road, Camberwell VIC 3124", "30/12/20 2:30pm-2:45pm", "Case shopped"),
    @("Cape Schank", "National Golf Club  The Cups Drive, Cape Schanck VIC 3939", "30/12/20 11.40am-1.40pm", "Case attended course"),
    @("Cheltenham", "Aldi Cheltenham  280/282 Bay Road, Cheltenham VIC 3192", "29/12/20 2:00pm-2:30pm", "Case shopped in store"),
    @("Cheltenham", "Angus and Cootes Jeweller  Southland Shopping Centre, 2096/1239 Nepean Hwy, Cheltenham VIC 3192", "28/12/2020 2:30pm-2:50pm", "Case shopped in store"),
    @("Cheltenham", "Honey Birdette Southland  Shop 2209/1239, Southland Shopping Centre, Cheltenham VIC 3192", "22/12/2020 3:50pm-4:05pm", "Case shopped in store"),
    @("Cheltenham", "Mecca Southland  Shop 2011/2013, Southland Shopping Centre, Cheltenham VIC 3192", "22/12/2020 3:30pm-3:50pm", "Case shopped in store"),
    @("Forest Hill", "Forest Hill Chase Shopping Centre 270 Canterbury Rd, Forest Hill VIC 3131", "28/12/20 12:00pm-2:00pm", "1210hrs Food court 30min; 1250hrs TKMaxx 15min; 1310hrs Target 20min; 1340hrs Woolworths 15min"),
    @("Fountain Gate Shopping Centre", "Kmart, Big W, Target, Millers, King of Gifts, Lo Costa  25-55 Overland Drive, Narre Warren VIC 3805", "26/12/20 9:00am-11:00am", ""),
    @("Glen Waverley", "Mocha Jos  87 Kingsway, Glen Waverley VIC 3150", "28/12/20 1:30pm-1:45pm", ""),
    @("Hallam", "Coles Hallam  2 Princes Domain Drive, Hallam VIC 3803", "30/12/20 6:15am-6:30am", "Case shopped in store"),
    @("Lakes Entrance", "Blue Riviera Hire Boats  Marine Parade, Lakes Entrance VIC 3909", "29/12/20 11:15am-12:15pm", "Case hired a boat"),
    @("Lakes Entrance", "Central Hotel Lakes Entrance  321 Esplanade, Lakes Entrance VIC 3909", "30/12/20 5:00pm-6:30pm", "Case attended outside premises"),
    @("Lakes Entrance", "Darcey Annas Beach Cafe Kiosk Gift Shop Gallery  426 Main Beach Walk Surf Life Saving, Lakes Entrance VIC 3909", "30/12/20 11:15am-11:20am", "Case picked up takeaway"),
    @("Lakes Entrance", "Woolworths Lakes Entrance 371 Esplanade, Lakes Entrance VIC 3909", "30/12/20 6:00pm-6:15pm", "Case shopped in store"),
    @("Mentone", "Bunnings Mentone  23-27 Nepean Hwy, Mentone VIC 3194", "29/12/20 07:30am-08:00am", "Case shopped in store"),
    @("Mentone", "Bunnings Mentone  23-27 Nepean Hwy, Mentone VIC 3194", "31/12/20 08:00am-08:30am", "Case shopped in store"),
    @("Mentone", "Mentone/Parkdale Beach", "27/12/20 10:00am-4:30pm", ""),
    @("Mentone", "Woolworths Mentone  105-111 Balcombe Road, Mentone VIC 3194", "23/12/20 2:45pm-3:05pm", "Case shopped in store"),
    @("Moorabbin", "COSTCO Moorabbin  8 Chifley Drive, Moorabbin Airport VIC 3194", "30/12/20 10:45am- 12:15pm", "Case shopped in store"),
    @("Moorabbin", "COSTCO Moorabbin  8 Chifley Drive, Moorabbin Airport VIC 3194", "30/12/20 4:00m- 5:50pm", "Case shopped in store"),
    @("Mordialloc", "Woodlands Golf Club  109 White Street Mordialloc VIC 3195", "23/12/20 8:00am-2:00pm", "Case attended course"),
    @("Mordialloc", "Woodlands Golf Club  109 White Street Mordialloc VIC 3195", "28/12/20 12:00pm-6:00pm", "Case attended course"),
    @("Mount Waverley", "Ritchies IGA  283 Stephensons Road, Mount Waverley VIC 3149", "30/12/20 2:00pm-2:30pm", "Case shopped for half an hour"),
    @("Oakleigh", "Bunnings Oakleigh  1041 Centre Road, Oakleigh South", "30/12/20 11:00am-11:30am", "Case shopped for 30 minutes"),
    @("Oakleigh", "Katialo restaurant  8 Eaton Mall, Oakleigh VIC 3166", "28/12/20 7:00pm-7:10pm", ""),
    @("Wonthaggi", "Wonthaggi Plaza Shopping centre  2 Biggs Drive, Wonthaggi VIC 3995", "28/12/20 1:30pm-2.30pm", "Kmart- shopped for 15 mins")
)

$startRow = 2
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $vals = $data[$i]
    $ws.Cells.Item($row, 1).Value = $vals[0]
    $ws.Cells.Item($row, 2).Value = $vals[1]
    $ws.Cells.Item($row, 3).Value = $vals[2]
    $ws.Cells.Item($row, 4).Value = $vals[3]
}
